# Apply the edit described by the commit:
#  "hadde skrivefeil i nøkkelord, rettet opp"  (had a typo in keyword, fixed it)
#
# 1. Insert a brand-new worksheet "programvare" right after "Statistikk"
#    (i.e. before " BIM"), populated with a course/learning-outcome table.
# 2. On the "Statistikk" sheet: fix the "programmvare" -> "programvare" typo
#    in A10, bump its hit counts (C10, E10: 0 -> 5), and update the total
#    hits in G2 (53 -> 58) to reflect the 5 new hits.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "programvare" worksheet, positioned before " BIM".
# ---------------------------------------------------------------------------
$bimSheet = $wb.Worksheets.Item(" BIM")
$ws = $wb.Worksheets.Add($bimSheet)
$ws.Name = "programvare"

# Column widths: 15 / 20 / 50 repeated for the three Emnekode/Emnenavn/LUT
# triples (columns A-C, D-F, G-I). COM ColumnWidth is offset by 5/6 from the
# width value persisted into the XML, so subtract that to land on the exact
# target widths of 15/20/50.
$offset = 5.0 / 6.0
$widths = @(15, 20, 50, 15, 20, 50, 15, 20, 50)
for ($i = 0; $i -lt $widths.Length; $i++) {
    $ws.Columns.Item($i + 1).ColumnWidth = $widths[$i] - $offset
}

# Header rows (LUK:/LUF:/LUG: and Emnekode:/Emnenavn:/Læringsutbytte triples)
$ws.Range("A1").Value = "LUK:"
$ws.Range("D1").Value = "LUF:"
$ws.Range("G1").Value = "LUG:"

$ws.Range("A2").Value = "Emnekode:"
$ws.Range("B2").Value = "Emnenavn:"
$ws.Range("C2").Value = "Læringsutbytte"
$ws.Range("D2").Value = "Emnekode:"
$ws.Range("E2").Value = "Emnenavn:"
$ws.Range("F2").Value = "Læringsutbytte"
$ws.Range("G2").Value = "Emnekode:"
$ws.Range("H2").Value = "Emnenavn:"
$ws.Range("I2").Value = "Læringsutbytte"

# Course data rows (only the D/E/F "LUF" columns are populated; A-C and G-I
# stay empty on this sheet).
$ws.Range("D3").Value = "BYVE3401"
$ws.Range("E3").Value = "Areal- og transportplanlegging"
$ws.Range("F3").Value = "Studenten kan utarbeide en reguleringsplan med tilhørende reguleringsbestemmelser og planbeskrivelse i samsvar med Miljøverndepartementets veiledning for reguleringsplaner samt overordnede føringer prinsippene for utforming av reguleringsplaner ved bruk av egnet programvare NovaPoint Areal Focus Arealplanlegging eller tilsvarende med tilhørende tekniske planer for Veg VA plantegning lengdeprofiler og tverrprofiler utføre konsekvensanalyser for områdereguleringsplan og ROS-analyse for detaljreguleringsplan utføre grunnleggende trafikktekniske beregninger og analyser"

$ws.Range("D4").Value = "EMPE1500"
$ws.Range("E4").Value = "Fysikk"
$ws.Range("F4").Value = "Studenten kan identifisere krefter og beregne kraftmomenter anvende Newtons 2 lov og spinnsatsen på konkrete fysiske problemer beskrive bevegelse matematisk blant annet ved hjelp av egnet programvare løse likevektproblemer for stive legemer"

$ws.Range("D5").Value = "EMPE2500"
$ws.Range("E5").Value = "Bygningssimulering"
$ws.Range("F5").Value = "Studenten kan utføre grunnleggende beregninger av varmetransport U-verdier kuldebro infiltrasjon og av effekt- og energibehovs med enkle formelverk håndtere bygningsinformasjonsmodeller BIM deriblant overføre data til programvare for simulering og miljøvurdering utføre dynamisk modellering av bygnings- og klimatekniske systemer for optimalt inneklima effekt- og energibehov ved bruk av simuleringsprogrammer som SIMIEN TEK-sjekk eller tilsvarende utføre energimerking av bygg vurdere inneklima termisk komfort og dagslysforhold utfra beregningene"

$ws.Range("D6").Value = "EMTS2600"
$ws.Range("E6").Value = "Inneklima og måleteknikk"
$ws.Range("F6").Value = "Studenten kan vurdere usikkerhet i alle typer målinger av inneklimaparametere og sette opp et usikkerhetsbudsjett håndtere spørreundersøkelser om inneklima ved hjelp av «Ørebroskjemaet» og tolke resultatet beregne nødvendige luftmengder ut ifra massebalanser og reaksjonskinetikk utføre målinger av inneklimaparametere som luftskifte luftkvalitet termiske akustiske og aktiniske forhold inkludert radon og sammenlikne dem med myndighetskrav vurdere materialbruken med hensyn på inneklimakvalitet og miljøbelastning foreta en mikrobiologisk analyse av en bygning spesielt med hensyn på muggsopp bruke Mollierediagram for å beregne duggpunkt og andre termodynamiske data for fuktig luft anvende programvare for inneklimasimuleringer designe for optimalt vedlikehold for å unngå Legionellavekst i varmtvannssystemer og kjøletårn designe våtrom"

$ws.Range("D7").Value = "EMVE3500"
$ws.Range("E7").Value = "Varme, ventilasjon og sanitærteknikk"
$ws.Range("F7").Value = "Studenten kan velge energikilderenergiforsyning som tilfredsstiller myndighetskrav gjennomføre energi- og inneklimaberegninger med relevant programvare utarbeide kravspesifikasjon for oppvarmingssystemer prosjekteredimensjonere energieffektive vannbårne oppvarmingsanlegg prosjekteredimensjonere energieffektive kjølesystemer utarbeide kravspesifikasjon for ventilasjonssystemer prosjekteredimensjonere energieffektive ventilasjonsanlegg herunder aggregat og kanalnett prosjekteredimensjonere ventilasjonsløsninger på rom nivå som gir akseptabelt inneklima med hensyn på temperatur trekk luftkvalitet og lyd prosjekteredimensjonere sanitærtekniske installasjoner innomhus vannforsyning og avløp prosjekteredimensjonere varmtvannsforsyningsanlegg"

# Whole-sheet formatting: vertical-top + wrap-text, matching the other
# course-listing sheets in this workbook.
$used = $ws.Range("A1:I7")
$used.VerticalAlignment = -4160   # xlTop
$used.WrapText = $true

# Page margins (inches 0.75/0.75/1/1/0.5/0.5 == points 54/54/72/72/36/36),
# matching the other course-listing sheets.
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------------
# 2. Fix the typo + counts on the "Statistikk" sheet.
# ---------------------------------------------------------------------------
$stat = $wb.Worksheets.Item("Statistikk")
$stat.Range("A10").Value = "programvare"
$stat.Range("C10").Value = 5
$stat.Range("E10").Value = 5
$stat.Range("G2").Value = 58

# Keep "Statistikk" as the active/selected tab, since adding a sheet makes
# the new sheet active by default and the source diff does not touch
# bookViews/activeTab.
$stat.Activate()
